$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 743.0769
$ws.Range("I19").Value = 541.4
$ws.Range("J19").Value = 1018.0909
$ws.Range("K19").Value = 541.4
$ws.Range("L19").Value = 1018.0909
$ws.Range("M19").Value = -366.4
$ws.Range("N19").Value = -1368.0909
$ws.Range("H113").Value = 93145.82000000001
$ws.Range("I113").Value = 201981
$ws.Range("J113").Value = 2449.8333
$ws.Range("K113").Value = 201981
$ws.Range("L113").Value = 2449.8333
$ws.Range("M113").Value = -198727
$ws.Range("N113").Value = -8957.8333
$ws.Range("H125").Value = 18686118
$ws.Range("I125").Value = 800
$ws.Range("J125").Value = 28028778
$ws.Range("K125").Value = 7200
$ws.Range("L125").Value = 252259002
$ws.Range("M125").Value = -4740
$ws.Range("N125").Value = -252263922
$ws.Range("H128").Value = 70000
$ws.Range("J128").Value = 70000
$ws.Range("L128").Value = 70000
$ws.Range("N128").Value = -79960
$ws.Range("H135").Value = 1114.6
$ws.Range("I135").Value = 1021.9286
$ws.Range("J135").Value = 2412
$ws.Range("K135").Value = 9197.357399999999
$ws.Range("L135").Value = 21708
$ws.Range("M135").Value = -6662.357399999999
$ws.Range("N135").Value = -26778
$ws.Range("H137").Value = 37038468
$ws.Range("I137").Value = 45455670
$ws.Range("K137").Value = 136367010
$ws.Range("M137").Value = -136364460

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15998.795
$ws.Range("I32").Value = 2378.446
$ws.Range("J32").Value = 126664.125
$ws.Range("K32").Value = 2378.446
$ws.Range("L32").Value = 126664.125
$ws.Range("M32").Value = -2091.446
$ws.Range("N32").Value = -127238.125
$ws.Range("H61").Value = 1823.3334
$ws.Range("I61").Value = 1182.3158
$ws.Range("J61").Value = 4259.2
$ws.Range("K61").Value = 1182.3158
$ws.Range("L61").Value = 4259.2
$ws.Range("M61").Value = -970.3158000000001
$ws.Range("N61").Value = -4683.2
$ws.Range("H136").Value = 1823.3334
$ws.Range("I136").Value = 1182.3158
$ws.Range("J136").Value = 4259.2
$ws.Range("K136").Value = 3546.9474
$ws.Range("L136").Value = 12777.6
$ws.Range("M136").Value = -996.9474
$ws.Range("N136").Value = -17877.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1712.0264
$ws.Range("I20").Value = 1783.2593
$ws.Range("J20").Value = 1537.1818
$ws.Range("K20").Value = 1783.2593
$ws.Range("L20").Value = 1537.1818
$ws.Range("M20").Value = -1536.2593
$ws.Range("N20").Value = -2031.1818
$ws.Range("H134").Value = 16951428
$ws.Range("I134").Value = 26317320
$ws.Range("J134").Value = 3625.8096
$ws.Range("K134").Value = 78951960
$ws.Range("L134").Value = 10877.4288
$ws.Range("M134").Value = -78949425
$ws.Range("N134").Value = -15947.4288

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2079.2144
$ws.Range("I31").Value = 1124.3529
$ws.Range("J31").Value = 3554.9092
$ws.Range("K31").Value = 1124.3529
$ws.Range("L31").Value = 3554.9092
$ws.Range("M31").Value = -829.3529000000001
$ws.Range("N31").Value = -4144.9092
$ws.Range("H34").Value = 2079.2144
$ws.Range("I34").Value = 1124.3529
$ws.Range("J34").Value = 3554.9092
$ws.Range("K34").Value = 1124.3529
$ws.Range("L34").Value = 3554.9092
$ws.Range("M34").Value = -922.3529000000001
$ws.Range("N34").Value = -3958.9092
$ws.Range("H58").Value = 2041.8286
$ws.Range("I58").Value = 576.125
$ws.Range("J58").Value = 3276.1052
$ws.Range("K58").Value = 576.125
$ws.Range("L58").Value = 3276.1052
$ws.Range("M58").Value = -373.125
$ws.Range("N58").Value = -3682.1052
$ws.Range("H99").Value = 8930288
$ws.Range("I99").Value = 20834500
$ws.Range("J99").Value = 2128.5
$ws.Range("K99").Value = 20834500
$ws.Range("L99").Value = 2128.5
$ws.Range("M99").Value = -20833002
$ws.Range("N99").Value = -5124.5
$ws.Range("H105").Value = 993.2143
$ws.Range("I105").Value = 915.8461
$ws.Range("J105").Value = 1999
$ws.Range("K105").Value = 915.8461
$ws.Range("L105").Value = 1999
$ws.Range("M105").Value = 831.1539
$ws.Range("N105").Value = -5493
$ws.Range("H126").Value = 8930288
$ws.Range("I126").Value = 20834500
$ws.Range("J126").Value = 2128.5
$ws.Range("K126").Value = 62503500
$ws.Range("L126").Value = 6385.5
$ws.Range("M126").Value = -62501030
$ws.Range("N126").Value = -11325.5
$ws.Range("H132").Value = 1993.1897
$ws.Range("I132").Value = 1453.7561
$ws.Range("J132").Value = 3294.1765
$ws.Range("K132").Value = 4361.2683
$ws.Range("L132").Value = 9882.529500000001
$ws.Range("M132").Value = -1831.2683
$ws.Range("N132").Value = -14942.5295
$ws.Range("H134").Value = 2070.8728
$ws.Range("I134").Value = 1257.1904
$ws.Range("J134").Value = 4699.6924
$ws.Range("K134").Value = 3771.5712
$ws.Range("L134").Value = 14099.0772
$ws.Range("M134").Value = -1236.5712
$ws.Range("N134").Value = -19169.0772
$ws.Range("H136").Value = 2041.8286
$ws.Range("I136").Value = 576.125
$ws.Range("J136").Value = 3276.1052
$ws.Range("K136").Value = 1728.375
$ws.Range("L136").Value = 9828.3156
$ws.Range("M136").Value = 821.625
$ws.Range("N136").Value = -14928.3156

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 17242358
$ws.Range("I113").Value = 672.3
$ws.Range("J113").Value = 26316928
$ws.Range("K113").Value = 2016.9
$ws.Range("L113").Value = 78950784
$ws.Range("M113").Value = 153.1000000000001
$ws.Range("N113").Value = -78955124
$ws.Range("H120").Value = 16250
$ws.Range("I120").Value = 3000
$ws.Range("K120").Value = 9000
$ws.Range("M120").Value = -4162

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2768.8333
$ws.Range("I113").Value = 2450
$ws.Range("J113").Value = 2928.25
$ws.Range("K113").Value = 2450
$ws.Range("L113").Value = 2928.25
$ws.Range("M113").Value = -280
$ws.Range("N113").Value = -7268.25
$ws.Range("H122").Value = 1112580.8
$ws.Range("I122").Value = 1390276
$ws.Range("K122").Value = 4170828
$ws.Range("M122").Value = -4168378
$ws.Range("H126").Value = 2805.55
$ws.Range("I126").Value = 2050
$ws.Range("K126").Value = 6150
$ws.Range("M126").Value = -3680
$ws.Range("H132").Value = 3059.5122
$ws.Range("I132").Value = 2836.6875
$ws.Range("J132").Value = 3851.7778
$ws.Range("K132").Value = 8510.0625
$ws.Range("L132").Value = 11555.3334
$ws.Range("M132").Value = -5980.0625
$ws.Range("N132").Value = -16615.3334
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 590
$ws.Range("I16").Value = 623.75
$ws.Range("K16").Value = 623.75
$ws.Range("M16").Value = -453.75
$ws.Range("H18").Value = 2168.3333
$ws.Range("I18").Value = 505
$ws.Range("K18").Value = 505
$ws.Range("M18").Value = -333
$ws.Range("H82").Value = 1146.0834
$ws.Range("J82").Value = 1161.4445
$ws.Range("L82").Value = 1161.4445
$ws.Range("N82").Value = -1883.4445
$ws.Range("H85").Value = 1146.0834
$ws.Range("J85").Value = 1161.4445
$ws.Range("L85").Value = 1161.4445
$ws.Range("N85").Value = -3657.4445
$ws.Range("H122").Value = 3768.158
$ws.Range("I122").Value = 1750
$ws.Range("J122").Value = 4005.5881
$ws.Range("K122").Value = 5250
$ws.Range("L122").Value = 12016.7643
$ws.Range("M122").Value = -2800
$ws.Range("N122").Value = -16916.7643
$ws.Range("H136").Value = 4107.3257
$ws.Range("I136").Value = 2460.6667
$ws.Range("J136").Value = 12575.857
$ws.Range("K136").Value = 7382.000100000001
$ws.Range("L136").Value = 37727.571
$ws.Range("M136").Value = -4832.000100000001
$ws.Range("N136").Value = -42827.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 273.15384
$ws.Range("I113").Value = 266.75
$ws.Range("J113").Value = 350
$ws.Range("K113").Value = 800.25
$ws.Range("L113").Value = 1050
$ws.Range("M113").Value = 1369.75
$ws.Range("N113").Value = -5390
$ws.Range("H122").Value = 68920.53
$ws.Range("I122").Value = 85650.664
$ws.Range("K122").Value = 256951.992
$ws.Range("M122").Value = -254501.992
$ws.Range("H132").Value = 11631295
$ws.Range("I132").Value = 16670250
$ws.Range("J132").Value = 2937.3845
$ws.Range("K132").Value = 50010750
$ws.Range("L132").Value = 8812.1535
$ws.Range("M132").Value = -50008220
$ws.Range("N132").Value = -13872.1535
$ws.Range("H136").Value = 10449237
$ws.Range("I136").Value = 11529681
$ws.Range("J136").Value = 4947.3335
$ws.Range("K136").Value = 34589043
$ws.Range("L136").Value = 14842.0005
$ws.Range("M136").Value = -34586493
$ws.Range("N136").Value = -19942.0005

